# Download Feature & File Update
# Mark several task-tracking checkboxes as "done" by writing the Wingdings
# checkmark glyph (U+00FC, displays as a tick in the Wingdings font) into
# the relevant cells and switching their font to Wingdings (matching the
# workbook's existing "checked" cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkMark = [char]0x00FC

$cells = @("J11", "J13", "J18", "G19", "H23", "J23", "J30")

foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Wingdings"
    $cell.Value = $checkMark
}
